$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1221.5714
$ws.Range("I32").Value = 1650
$ws.Range("J32").Value = 1150.1666
$ws.Range("K32").Value = 1650
$ws.Range("L32").Value = 1150.1666
$ws.Range("M32").Value = -1324
$ws.Range("N32").Value = -1802.1666
$ws.Range("H137").Value = 5406216
$ws.Range("I137").Value = 610.64
$ws.Range("K137").Value = 1831.92
$ws.Range("M137").Value = 718.0799999999999
$ws.Range("H139").Value = 200000
$ws.Range("J139").Value = 200000
$ws.Range("L139").Value = 200000
$ws.Range("N139").Value = -210280
$ws.Range("H140").Value = 98333.336
$ws.Range("J140").Value = 98333.336
$ws.Range("L140").Value = 98333.336
$ws.Range("N140").Value = -108693.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10665.972
$ws.Range("I32").Value = 11620.72
$ws.Range("J32").Value = 8279.1
$ws.Range("K32").Value = 11620.72
$ws.Range("L32").Value = 8279.1
$ws.Range("M32").Value = -11333.72
$ws.Range("N32").Value = -8853.1
$ws.Range("H61").Value = 10871121
$ws.Range("I61").Value = 12501591
$ws.Range("J61").Value = 1319.1666
$ws.Range("K61").Value = 12501591
$ws.Range("L61").Value = 1319.1666
$ws.Range("M61").Value = -12501379
$ws.Range("N61").Value = -1743.1666
$ws.Range("H110").Value = 1100.5834
$ws.Range("I110").Value = 851.7
$ws.Range("J110").Value = 2345
$ws.Range("K110").Value = 851.7
$ws.Range("L110").Value = 2345
$ws.Range("M110").Value = 1193.3
$ws.Range("N110").Value = -6435
$ws.Range("H122").Value = 7458.95
$ws.Range("I122").Value = 8814.933999999999
$ws.Range("J122").Value = 3391
$ws.Range("K122").Value = 26444.802
$ws.Range("L122").Value = 10173
$ws.Range("M122").Value = -23994.802
$ws.Range("N122").Value = -15073
$ws.Range("H132").Value = 3379815
$ws.Range("J132").Value = 1589
$ws.Range("L132").Value = 4767
$ws.Range("N132").Value = -9827
$ws.Range("H136").Value = 10871121
$ws.Range("I136").Value = 12501591
$ws.Range("J136").Value = 1319.1666
$ws.Range("K136").Value = 37504773
$ws.Range("L136").Value = 3957.4998
$ws.Range("M136").Value = -37502223
$ws.Range("N136").Value = -9057.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 17262.2
$ws.Range("I82").Value = 14128.5
$ws.Range("J82").Value = 19351.334
$ws.Range("K82").Value = 14128.5
$ws.Range("L82").Value = 19351.334
$ws.Range("M82").Value = -13745.5
$ws.Range("N82").Value = -20117.334
$ws.Range("H85").Value = 17262.2
$ws.Range("I85").Value = 14128.5
$ws.Range("J85").Value = 19351.334
$ws.Range("K85").Value = 14128.5
$ws.Range("L85").Value = 19351.334
$ws.Range("M85").Value = -12802.5
$ws.Range("N85").Value = -22003.334
$ws.Range("H105").Value = 1964.9104
$ws.Range("I105").Value = 986.0244
$ws.Range("J105").Value = 3508.5386
$ws.Range("K105").Value = 986.0244
$ws.Range("L105").Value = 3508.5386
$ws.Range("M105").Value = 760.9756
$ws.Range("N105").Value = -7002.5386
$ws.Range("H107").Value = 1041.0435
$ws.Range("I107").Value = 879.05884
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 879.05884
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 1040.94116
$ws.Range("N107").Value = -5340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 500511
$ws.Range("I16").Value = 500511
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 500511
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -500224
$ws.Range("N16").ClearContents()
$ws.Range("H20").Value = 38900
$ws.Range("J20").Value = 38900
$ws.Range("L20").Value = 38900
$ws.Range("N20").Value = -39372
$ws.Range("H30").Value = 38900
$ws.Range("J30").Value = 38900
$ws.Range("L30").Value = 38900
$ws.Range("N30").Value = -39082
$ws.Range("H31").Value = 6671175
$ws.Range("I31").Value = 4958.057
$ws.Range("K31").Value = 4958.057
$ws.Range("M31").Value = -4663.057
$ws.Range("H34").Value = 6671175
$ws.Range("I34").Value = 4958.057
$ws.Range("K34").Value = 4958.057
$ws.Range("M34").Value = -4756.057
$ws.Range("H62").Value = 2235.3333
$ws.Range("I62").Value = 2192
$ws.Range("J62").Value = 2452
$ws.Range("K62").Value = 2192
$ws.Range("L62").Value = 2452
$ws.Range("M62").Value = -1568
$ws.Range("N62").Value = -3700
$ws.Range("H65").Value = 2235.3333
$ws.Range("I65").Value = 2192
$ws.Range("J65").Value = 2452
$ws.Range("K65").Value = 10960
$ws.Range("L65").Value = 12260
$ws.Range("M65").Value = -7840
$ws.Range("N65").Value = -18500
$ws.Range("H113").Value = 500511
$ws.Range("I113").Value = 500511
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 500511
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -498341
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2703.25
$ws.Range("I122").Value = 2703.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8109.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5659.75
$ws.Range("N122").ClearContents()
$ws.Range("H128").Value = 38900
$ws.Range("J128").Value = 38900
$ws.Range("L128").Value = 38900
$ws.Range("N128").Value = -48860
$ws.Range("H140").Value = 24335.8
$ws.Range("J140").Value = 24335.8
$ws.Range("L140").Value = 24335.8
$ws.Range("N140").Value = -34695.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 262
$ws.Range("I51").Value = 262
$ws.Range("K51").Value = 786
$ws.Range("M51").Value = -326
$ws.Range("H63").Value = 9744.75
$ws.Range("I63").Value = 9699
$ws.Range("J63").Value = 9760
$ws.Range("K63").Value = 29097
$ws.Range("L63").Value = 29280
$ws.Range("M63").Value = -28348
$ws.Range("N63").Value = -30778
$ws.Range("H66").Value = 9744.75
$ws.Range("I66").Value = 9699
$ws.Range("J66").Value = 9760
$ws.Range("K66").Value = 87291
$ws.Range("L66").Value = 87840
$ws.Range("M66").Value = -83547
$ws.Range("N66").Value = -95328

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H132").Value = 3175.0833
$ws.Range("I132").Value = 2245.6758
$ws.Range("J132").Value = 6301.273
$ws.Range("K132").Value = 6737.0274
$ws.Range("L132").Value = 18903.819
$ws.Range("M132").Value = -4207.0274
$ws.Range("N132").Value = -23963.819
$ws.Range("H138").Value = 59800
$ws.Range("J138").Value = 59800
$ws.Range("L138").Value = 59800
$ws.Range("N138").Value = -70080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6636
$ws.Range("I40").Value = 6977.6665
$ws.Range("J40").Value = 6021
$ws.Range("K40").Value = 6977.6665
$ws.Range("L40").Value = 6021
$ws.Range("M40").Value = -6841.6665
$ws.Range("N40").Value = -6293
$ws.Range("H46").Value = 787.625
$ws.Range("I46").Value = 544.2
$ws.Range("J46").Value = 1193.3334
$ws.Range("K46").Value = 544.2
$ws.Range("L46").Value = 1193.3334
$ws.Range("M46").Value = -356.2
$ws.Range("N46").Value = -1569.3334
$ws.Range("H61").Value = 1471.7142
$ws.Range("I61").Value = 1262.375
$ws.Range("K61").Value = 1262.375
$ws.Range("M61").Value = -1060.375
$ws.Range("H113").Value = 1471.7142
$ws.Range("I113").Value = 1262.375
$ws.Range("K113").Value = 1262.375
$ws.Range("M113").Value = 907.625
$ws.Range("H132").Value = 7048026
$ws.Range("I132").Value = 3290.549
$ws.Range("J132").Value = 25012102
$ws.Range("K132").Value = 9871.647000000001
$ws.Range("L132").Value = 75036306
$ws.Range("M132").Value = -7341.647000000001
$ws.Range("N132").Value = -75041366
$ws.Range("H139").Value = 57090
$ws.Range("J139").Value = 57090
$ws.Range("L139").Value = 57090
$ws.Range("N139").Value = -67370
